$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert a new column at F ("GRANITO") ---
$ws.Range("F1").EntireColumn.Insert()

# --- Insert three new columns after the (shifted) N column, i.e. at P:R
#     ("NO RESURTIBLES", "PANELES PVC", "PANELES PU") ---
$ws.Range("P1:R1").EntireColumn.Insert()

# --- Column widths (character units stored in xlsx = ColumnWidth + 5/6) ---
$offset = 0.8333333333333334
$ws.Range("F1").ColumnWidth = 13 - $offset
$ws.Range("G1").ColumnWidth = 15 - $offset
$ws.Range("H1").ColumnWidth = 14 - $offset
$ws.Range("I1").ColumnWidth = 13 - $offset
$ws.Range("J1").ColumnWidth = 9 - $offset
$ws.Range("K1").ColumnWidth = 25 - $offset
$ws.Range("L1").ColumnWidth = 24 - $offset
$ws.Range("M1").ColumnWidth = 17 - $offset
$ws.Range("N1").ColumnWidth = 26 - $offset
$ws.Range("O1").ColumnWidth = 17 - $offset
$ws.Range("P1").ColumnWidth = 20 - $offset
$ws.Range("Q1").ColumnWidth = 17 - $offset
$ws.Range("R1").ColumnWidth = 16 - $offset

# --- Header row text for the new columns ---
$ws.Range("F1").Value = "GRANITO"
$ws.Range("P1").Value = "NO RESURTIBLES"
$ws.Range("Q1").Value = "PANELES PVC"
$ws.Range("R1").Value = "PANELES PU"

# --- Data rows 2-21: the new column F and the new columns P:R are 0
#     (column O already holds the values shifted over from the old N column) ---
$ws.Range("F2:F21").Value = 0
$ws.Range("P2:R21").Value = 0

# --- Footer row 22: the new column F and the new columns P:R are "0 de 20"
#     (column O already holds the "0 de 20" text shifted over from old N) ---
$ws.Range("F22").Value = "0 de 20"
$ws.Range("P22:R22").Value = "0 de 20"
